$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3845.9285
$ws.Range("I74").Value = 3510
$ws.Range("J74").Value = 4181.857
$ws.Range("K74").Value = 3510
$ws.Range("L74").Value = 4181.857
$ws.Range("M74").Value = -2574
$ws.Range("N74").Value = -6053.857

$ws.Range("H77").Value = 3845.9285
$ws.Range("I77").Value = 3510
$ws.Range("J77").Value = 4181.857
$ws.Range("K77").Value = 17550
$ws.Range("L77").Value = 20909.285
$ws.Range("M77").Value = -12870
$ws.Range("N77").Value = -30269.285

$ws.Range("H113").Value = 5001.6665
$ws.Range("I113").Value = 5001.6665
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 5001.6665
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -1747.6665

$ws.Range("H129").Value = 1074.6182
$ws.Range("I129").Value = 552.5833
$ws.Range("J129").Value = 1220.3024
$ws.Range("K129").Value = 1657.7499
$ws.Range("L129").Value = 3660.9072
$ws.Range("M129").Value = 3342.2501
$ws.Range("N129").Value = -13660.9072

$ws.Range("H138").Value = 1975.1266
$ws.Range("I138").Value = 1296.1333
$ws.Range("J138").Value = 2873.7942
$ws.Range("K138").Value = 3888.3999
$ws.Range("L138").Value = 8621.382599999999
$ws.Range("M138").Value = 1251.6001
$ws.Range("N138").Value = -18901.3826

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1734.6
$ws.Range("I74").Value = 1447.6666
$ws.Range("J74").Value = 2308.4666
$ws.Range("K74").Value = 1447.6666
$ws.Range("L74").Value = 2308.4666
$ws.Range("M74").Value = -573.6666
$ws.Range("N74").Value = -4056.4666

$ws.Range("H77").Value = 1734.6
$ws.Range("I77").Value = 1447.6666
$ws.Range("J77").Value = 2308.4666
$ws.Range("K77").Value = 7238.333000000001
$ws.Range("L77").Value = 11542.333
$ws.Range("M77").Value = -2870.333000000001
$ws.Range("N77").Value = -20278.333

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 15287.125
$ws.Range("I82").Value = 7235.6665
$ws.Range("J82").Value = 20118
$ws.Range("K82").Value = 7235.6665
$ws.Range("L82").Value = 20118
$ws.Range("M82").Value = -6852.6665
$ws.Range("N82").Value = -20884

$ws.Range("H85").Value = 15287.125
$ws.Range("I85").Value = 7235.6665
$ws.Range("J85").Value = 20118
$ws.Range("K85").Value = 7235.6665
$ws.Range("L85").Value = 20118
$ws.Range("M85").Value = -5909.6665
$ws.Range("N85").Value = -22770

$ws.Range("H94").Value = 1605.8125
$ws.Range("I94").Value = 1399.3
$ws.Range("J94").Value = 1950
$ws.Range("K94").Value = 1399.3
$ws.Range("L94").Value = 1950
$ws.Range("M94").Value = -948.3
$ws.Range("N94").Value = -2852

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 27625.062
$ws.Range("J4").Value = 27800.066
$ws.Range("L4").Value = 27800.066
$ws.Range("N4").Value = -28024.066

$ws.Range("H43").Value = 288000
$ws.Range("J43").Value = 288000
$ws.Range("L43").Value = 288000
$ws.Range("N43").Value = -288368

$ws.Range("H92").Value = 80000
$ws.Range("J92").Value = 80000
$ws.Range("L92").Value = 80000
$ws.Range("N92").Value = -84992

$ws.Range("H101").Value = 288000
$ws.Range("J101").Value = 288000
$ws.Range("L101").Value = 288000
$ws.Range("N101").Value = -294490

$ws.Range("H132").Value = 2223274.8
$ws.Range("I132").Value = 941.95746
$ws.Range("J132").Value = 5953619
$ws.Range("K132").Value = 2825.87238
$ws.Range("L132").Value = 17860857
$ws.Range("M132").Value = -295.8723799999998
$ws.Range("N132").Value = -17865917

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 14001250
$ws.Range("I4").Value = 17500850
$ws.Range("J4").Value = 2850
$ws.Range("K4").Value = 52502550
$ws.Range("L4").Value = 8550
$ws.Range("M4").Value = -52502438
$ws.Range("N4").Value = -8774

$ws.Range("H5").Value = 917.9231
$ws.Range("I5").Value = 572.4
$ws.Range("K5").Value = 1717.2
$ws.Range("M5").Value = -1605.2

$ws.Range("H31").Value = 1558.5
$ws.Range("J31").Value = 1558.5
$ws.Range("L31").Value = 4675.5
$ws.Range("N31").Value = -5251.5

$ws.Range("H107").Value = 58823976
$ws.Range("I107").Value = 358.57144
$ws.Range("J107").Value = 100000504
$ws.Range("K107").Value = 1075.71432
$ws.Range("L107").Value = 300001512
$ws.Range("M107").Value = 844.28568
$ws.Range("N107").Value = -300005352

$ws.Range("H122").Value = 5412.591
$ws.Range("I122").Value = 535
$ws.Range("J122").Value = 10290.182
$ws.Range("K122").Value = 4815
$ws.Range("L122").Value = 92611.63800000001
$ws.Range("M122").Value = -2365
$ws.Range("N122").Value = -97511.63800000001

$ws.Range("H125").Value = 2670.625
$ws.Range("J125").Value = 2909.2856
$ws.Range("L125").Value = 8727.856800000001
$ws.Range("N125").Value = -18567.8568

$ws.Range("H132").Value = 2323.1
$ws.Range("I132").Value = 2096.4614
$ws.Range("J132").Value = 2432.2222
$ws.Range("K132").Value = 18868.1526
$ws.Range("L132").Value = 21889.9998
$ws.Range("M132").Value = -16338.1526
$ws.Range("N132").Value = -26949.9998

$ws.Range("H135").Value = 917.9231
$ws.Range("I135").Value = 572.4
$ws.Range("K135").Value = 5151.599999999999
$ws.Range("M135").Value = -2616.599999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 63676500
$ws.Range("J80").Value = 102249.75
$ws.Range("L80").Value = 102249.75
$ws.Range("N80").Value = -104245.75

$ws.Range("H83").Value = 63676500
$ws.Range("J83").Value = 102249.75
$ws.Range("L83").Value = 511248.75
$ws.Range("N83").Value = -521232.75

$ws.Range("H132").Value = 2048.0945
$ws.Range("I132").Value = 1817.8704
$ws.Range("J132").Value = 2669.7
$ws.Range("K132").Value = 5453.6112
$ws.Range("L132").Value = 8009.099999999999
$ws.Range("M132").Value = -2923.6112
$ws.Range("N132").Value = -13069.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 63334.668
$ws.Range("J2").Value = 63334.668
$ws.Range("L2").Value = 63334.668
$ws.Range("N2").Value = -63558.668

$ws.Range("H40").Value = 113287.336
$ws.Range("I40").Value = 144799.42
$ws.Range("J40").Value = 2995
$ws.Range("K40").Value = 144799.42
$ws.Range("L40").Value = 2995
$ws.Range("M40").Value = -144663.42
$ws.Range("N40").Value = -3267

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 39990
$ws.Range("J82").Value = 39990
$ws.Range("L82").Value = 39990
$ws.Range("N82").Value = -40756

$ws.Range("H85").Value = 39990
$ws.Range("J85").Value = 39990
$ws.Range("L85").Value = 39990
$ws.Range("N85").Value = -42642
